$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and a two-row reorder at 34/35)
$ws.Range("D2").Value = "36.279.82"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "1.929.27"
$ws.Range("E3").Value = "  -2.80%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'240.07"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").Value = "'0.602"
$ws.Range("E6").Value = "  -3.97%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'55.54"
$ws.Range("E8").Value = "  -6.14%  "
$ws.Range("E9").Value = "  -5.66%  "
$ws.Range("D10").Value = "'0.0826"
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("D11").Value = "'0.102"
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("D12").Value = "2.212.58"
$ws.Range("E12").Value = "  -2.87%  "
$ws.Range("D13").Value = "'0.790"
$ws.Range("E13").Value = "  -8.56%  "
$ws.Range("D14").Value = "'13.17"
$ws.Range("E14").Value = "  -5.86%  "
$ws.Range("D15").Value = "'20.57"
$ws.Range("D16").Value = "'5.06"
$ws.Range("E16").Value = "  -7.29%  "
$ws.Range("D17").Value = "1.929.42"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").Value = "36.217.11"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "'68.34"
$ws.Range("E19").Value = "  -3.08%  "
$ws.Range("D20").Value = "0.0₃0851"
$ws.Range("E20").Value = "  -2.46%  "
$ws.Range("D21").Value = "'225.20"
$ws.Range("E21").Value = "  -3.79%  "
$ws.Range("D22").Value = "'4.89"
$ws.Range("E22").Value = "  -7.97%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'2.30"
$ws.Range("E24").Value = "  -10.24%  "
$ws.Range("D25").Value = "'2.23"
$ws.Range("E25").Value = "  -2.69%  "
$ws.Range("D26").Value = "'9.00"
$ws.Range("E26").Value = "  -9.56%  "
$ws.Range("D27").Value = "'159.86"
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("D28").Value = "'0.129"
$ws.Range("E28").Value = "  -3.22%  "
$ws.Range("D29").Value = "'18.94"
$ws.Range("E29").Value = "  -4.64%  "
$ws.Range("E30").Value = "  -3.16%  "
$ws.Range("D31").Value = "'1.08"
$ws.Range("E31").Value = "  -8.10%  "
$ws.Range("D32").Value = "'4.47"
$ws.Range("E32").Value = "  -8.71%  "
$ws.Range("D33").Value = "'0.0614"
$ws.Range("E33").Value = "  -7.81%  "
$ws.Range("B34").Value = "BinanceUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'4.09"
$ws.Range("E35").Value = "  -7.27%  "
$ws.Range("E36").Value = "  -3.64%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  -6.64%  "
$ws.Range("D39").Value = "'2.93"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").Value = "'0.0955"
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").Value = "'2.85"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").Value = "'0.0207"
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("E43").Value = "  -8.57%  "
$ws.Range("D44").Value = "'15.32"
$ws.Range("E44").Value = "  -5.53%  "
$ws.Range("D45").Value = "1.324.34"
$ws.Range("E45").Value = "  -3.07%  "
$ws.Range("E46").Value = "  -8.28%  "
$ws.Range("D47").Value = "'6.99"
$ws.Range("E47").Value = "  -6.54%  "
$ws.Range("D48").Value = "'84.02"
$ws.Range("E48").Value = "  -8.97%  "
$ws.Range("D49").Value = "'2.81"
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("D50").Value = "2.104.94"
$ws.Range("E50").Value = "  -2.78%  "
$ws.Range("D51").Value = "'42.57"
$ws.Range("E51").Value = "  -5.86%  "
